$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "Model" label in E2 (ProtocolType column was left over from a copy/paste;
# only PlayerId/Id keep their Model type note in E3 going forward).
$null = $ws.Range("E2").ClearContents()

# Leave the cursor parked on F5 before saving, matching where editing left off.
$null = $ws.Range("F5").Select()
